# Apply update "04/01/2017 Update after GW call." to Gemalto-CloudGate-Issues_O.xlsx
$wb = $excel.ActiveWorkbook

$wsIssues = $wb.Worksheets.Item("Issue Tracking")
$wsLegend = $wb.Worksheets.Item("Legend")

# 1. Insert a new Status option "Under Test (Customer)" into the Legend sheet,
#    just before the existing "New firmware" entry (Legend!C7), pushing it down to C8.
$wsLegend.Range("C7").Insert(-4121)
$wsLegend.Range("C7").Value = "Under Test (Customer)"

# 2. Update the "Status" for issues 10 and 11 (rows 12 & 13) on the Issue Tracking sheet
#    from "Under Test (GW)" to the new "Under Test (Customer)" status.
$wsIssues.Range("G12").Value = "Under Test (Customer)"
$wsIssues.Range("G13").Value = "Under Test (Customer)"

# 3. Update the Resolution/Plan comments for customer issue #6 (row 8) with the latest
#    status from Jimmy, including the new 01/04 update.
$wsIssues.Range("E8").Value = "•Jimmy needs to have a teamviewer session to debug what is going on here.`n•11/28 Teamviewer ready for Jimmy.`n•11/29 Jimmy tested and saw a delay of 2 minutes when switching. He will investigate and see if there is any chance for improvement.`n•There is a problem with this setup being offline every day. Jimmy is loosing time as he cannot test without an active TV.`n•01/04 Jimmy is going to try to see the same issue in Belgium and if he does the TV setup will no longer be needed."
$wsIssues.Rows.Item(8).RowHeight = 195

# 4. Update the Resolution/Plan comments for customer issue #9 (row 11) with Brandon's
#    01/04 update.
$wsIssues.Range("E11").Value = "•12/22 Brandon to Verify if this is the case.`n•01/04 Brandon is waiting for his device to come back in order to test this."

# 5. Refresh the view state: Legend sheet scrolled to show the newly added row, and the
#    Issue Tracking sheet scrolled down with G8 selected (per the saved workbook state).
$wsLegend.Range("D10").Select()
$wsIssues.Activate()
$wsIssues.Range("G8").Select()
